$p = $ppt.ActivePresentation

# Add a new slide at the end (position 11) using the "Title and Content" layout,
# matching ppt/slideLayouts/slideLayout6.xml's family (Title Only) -- we start
# from Title Only so the Slide Number placeholder gets id=3, then paste in the
# Content Placeholder afterwards so it gets id=4 (matching original authoring order).
$s = $p.Slides.Add(11, 11)

# --- Title ---
$s.Shapes.Item(1).TextFrame.TextRange.Text = "What was left out?"

# --- Slide Number placeholder: copy from slide 10 so it keeps the same field
# id/name/placeholder-id numbering already used throughout this deck ---
$numSrc = $p.Slides.Item(10).Shapes.Item(2)
$numSrc.Copy()
$s.Shapes.Paste() | Out-Null
$numShape = $s.Shapes.Item($s.Shapes.Count)

# --- Content placeholder: copy an existing "Content Placeholder" shape so the
# new slide gets a real placeholder shape (idx=1) rather than a plain textbox ---
$contentSrc = $p.Slides.Item(5).Shapes.Item(2)
$contentSrc.Copy()
$s.Shapes.Paste() | Out-Null
$contentShape = $s.Shapes.Item($s.Shapes.Count)
$contentShape.Name = "Content Placeholder 3"

# Put the content placeholder ahead of the slide number placeholder in the
# shape order (but keep their already-assigned ids).
$contentShape.ZOrder(3)

# --- Set the bullet text for the content placeholder ---
$tr = $contentShape.TextFrame.TextRange
$tr.Text = "Abstract classes" + [char]13 + "virtual functions" + [char]13 + "Multiple inheritance"
$tr.Paragraphs(2, 1).IndentLevel = 2
